$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 885.9655
$ws.Range("J17").Value = 885.9655
$ws.Range("L17").Value = 2657.8965
$ws.Range("N17").Value = -2993.8965
$ws.Range("H32").Value = 1653.2222
$ws.Range("I32").Value = 4000
$ws.Range("J32").Value = 1359.875
$ws.Range("K32").Value = 4000
$ws.Range("L32").Value = 1359.875
$ws.Range("M32").Value = -3674
$ws.Range("N32").Value = -2011.875
$ws.Range("H40").Value = 2400
$ws.Range("I40").Value = 2666.6667
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 2666.6667
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -2491.6667
$ws.Range("N40").Value = -2350
$ws.Range("H98").Value = 1865.9584
$ws.Range("I98").Value = 1965.9048
$ws.Range("J98").Value = 1166.3334
$ws.Range("K98").Value = 1965.9048
$ws.Range("L98").Value = 1166.3334
$ws.Range("M98").Value = -467.9048
$ws.Range("N98").Value = -4162.3334
$ws.Range("H108").Value = 63999
$ws.Range("J108").Value = 63999
$ws.Range("L108").Value = 63999
$ws.Range("N108").Value = -71679
$ws.Range("H109").Value = 30684
$ws.Range("J109").Value = 30684
$ws.Range("L109").Value = 30684
$ws.Range("N109").Value = -33458
$ws.Range("H122").Value = 1865.9584
$ws.Range("I122").Value = 1965.9048
$ws.Range("J122").Value = 1166.3334
$ws.Range("K122").Value = 5897.7144
$ws.Range("L122").Value = 3499.0002
$ws.Range("M122").Value = -3447.7144
$ws.Range("N122").Value = -8399.0002
$ws.Range("H131").Value = 3197.8333
$ws.Range("J131").Value = 4586.4287
$ws.Range("L131").Value = 13759.2861
$ws.Range("N131").Value = -23839.2861
$ws.Range("H137").Value = 1313.7742
$ws.Range("I137").Value = 1231.8077
$ws.Range("J137").Value = 1740
$ws.Range("K137").Value = 3695.4231
$ws.Range("L137").Value = 5220
$ws.Range("M137").Value = -1145.4231
$ws.Range("N137").Value = -10320
$ws.Range("H140").Value = 51271.2
$ws.Range("J140").Value = 51271.2
$ws.Range("L140").Value = 51271.2
$ws.Range("N140").Value = -61631.2

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 38.75
$ws.Range("I4").Value = 38.75
$ws.Range("K4").Value = 38.75
$ws.Range("M4").Value = 77.25
$ws.Range("H5").Value = 222
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H32").Value = 3351.5942
$ws.Range("I32").Value = 2533.6167
$ws.Range("K32").Value = 2533.6167
$ws.Range("M32").Value = -2246.6167
$ws.Range("H61").Value = 4704.7856
$ws.Range("I61").Value = 3157.5715
$ws.Range("K61").Value = 3157.5715
$ws.Range("M61").Value = -2945.5715
$ws.Range("H92").Value = 32775
$ws.Range("J92").Value = 32775
$ws.Range("L92").Value = 32775
$ws.Range("N92").Value = -37767
$ws.Range("H132").Value = 2443.6086
$ws.Range("I132").Value = 2134.6667
$ws.Range("K132").Value = 6404.000100000001
$ws.Range("M132").Value = -3874.000100000001
$ws.Range("H136").Value = 4704.7856
$ws.Range("I136").Value = 3157.5715
$ws.Range("K136").Value = 9472.7145
$ws.Range("M136").Value = -6922.7145
$ws.Range("H139").Value = 34500
$ws.Range("J139").Value = 34500
$ws.Range("L139").Value = 34500
$ws.Range("N139").Value = -44780

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 222
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H99").Value = 1220.1052
$ws.Range("I99").Value = 1014.1539
$ws.Range("K99").Value = 1014.1539
$ws.Range("M99").Value = 483.8461
$ws.Range("H105").Value = 2194.1304
$ws.Range("I105").Value = 2165
$ws.Range("K105").Value = 2165
$ws.Range("M105").Value = -418
$ws.Range("H107").Value = 1686.7693
$ws.Range("I107").Value = 1686.7693
$ws.Range("K107").Value = 1686.7693
$ws.Range("M107").Value = 233.2307000000001
$ws.Range("H108").Value = 94979.5
$ws.Range("J108").Value = 94979.5
$ws.Range("L108").Value = 94979.5
$ws.Range("N108").Value = -102659.5
$ws.Range("H134").Value = 2462.4119
$ws.Range("I134").Value = 2567.4
$ws.Range("J134").Value = 1675
$ws.Range("K134").Value = 7702.200000000001
$ws.Range("L134").Value = 5025
$ws.Range("M134").Value = -5167.200000000001
$ws.Range("N134").Value = -10095

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 20000
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H94").Value = 1104.8
$ws.Range("J94").Value = 1362
$ws.Range("L94").Value = 1362
$ws.Range("N94").Value = -2264

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 159.77777
$ws.Range("I2").Value = 164.83333
$ws.Range("K2").Value = 988.9999799999999
$ws.Range("M2").Value = -875.9999799999999
$ws.Range("H4").Value = 400079.88
$ws.Range("I4").Value = 400079.88
$ws.Range("K4").Value = 1200239.64
$ws.Range("M4").Value = -1200127.64
$ws.Range("H37").Value = 100000
$ws.Range("J37").Value = 100000
$ws.Range("L37").Value = 300000
$ws.Range("N37").Value = -300224
$ws.Range("H98").Value = 986.75
$ws.Range("I98").Value = 398
$ws.Range("J98").Value = 1183
$ws.Range("K98").Value = 1194
$ws.Range("L98").Value = 3549
$ws.Range("M98").Value = 304
$ws.Range("N98").Value = -6545
$ws.Range("H122").Value = 968.5454999999999
$ws.Range("J122").Value = 1074.375
$ws.Range("L122").Value = 9669.375
$ws.Range("N122").Value = -14569.375

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 124.125
$ws.Range("I2").Value = 170.85715
$ws.Range("J2").Value = 87.77778000000001
$ws.Range("K2").Value = 170.85715
$ws.Range("L2").Value = 87.77778000000001
$ws.Range("M2").Value = -57.85714999999999
$ws.Range("N2").Value = -313.77778
$ws.Range("H122").Value = 2157
$ws.Range("I122").Value = 1719.8
$ws.Range("K122").Value = 5159.4
$ws.Range("M122").Value = -2709.4
$ws.Range("H132").Value = 3207651
$ws.Range("I132").Value = 5496148
$ws.Range("J132").Value = 3755.4
$ws.Range("K132").Value = 16488444
$ws.Range("L132").Value = 11266.2
$ws.Range("M132").Value = -16485914
$ws.Range("N132").Value = -16326.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2007.1428
$ws.Range("I22").Value = 2204.1667
$ws.Range("J22").Value = 1859.375
$ws.Range("K22").Value = 2204.1667
$ws.Range("L22").Value = 1859.375
$ws.Range("M22").Value = -1909.1667
$ws.Range("N22").Value = -2449.375
$ws.Range("H27").Value = 2007.1428
$ws.Range("I27").Value = 2204.1667
$ws.Range("J27").Value = 1859.375
$ws.Range("K27").Value = 2204.1667
$ws.Range("L27").Value = 1859.375
$ws.Range("M27").Value = -2097.1667
$ws.Range("N27").Value = -2073.375
$ws.Range("H32").Value = 6646.8
$ws.Range("I32").Value = 5841
$ws.Range("K32").Value = 5841
$ws.Range("M32").Value = -5524
$ws.Range("H40").Value = 3785
$ws.Range("I40").Value = 1542
$ws.Range("J40").Value = 15000
$ws.Range("K40").Value = 1542
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = -1406
$ws.Range("N40").Value = -15272
$ws.Range("H46").Value = 2465.4
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 2517.111
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 2517.111
$ws.Range("M46").Value = -1812
$ws.Range("N46").Value = -2893.111
$ws.Range("H82").Value = 2375.8333
$ws.Range("I82").Value = 1750.4
$ws.Range("J82").Value = 2822.5715
$ws.Range("K82").Value = 1750.4
$ws.Range("L82").Value = 2822.5715
$ws.Range("M82").Value = -1389.4
$ws.Range("N82").Value = -3544.5715
$ws.Range("H85").Value = 2375.8333
$ws.Range("I85").Value = 1750.4
$ws.Range("J85").Value = 2822.5715
$ws.Range("K85").Value = 1750.4
$ws.Range("L85").Value = 2822.5715
$ws.Range("M85").Value = -502.4000000000001
$ws.Range("N85").Value = -5318.5715
$ws.Range("H94").Value = 48329.5
$ws.Range("J94").Value = 48329.5
$ws.Range("L94").Value = 48329.5
$ws.Range("N94").Value = -49681.5
$ws.Range("H122").Value = 12499.75
$ws.Range("I122").Value = 9999.5
$ws.Range("K122").Value = 29998.5
$ws.Range("M122").Value = -27548.5
$ws.Range("H128").Value = 100429
$ws.Range("J128").Value = 100429
$ws.Range("L128").Value = 100429
$ws.Range("N128").Value = -110389

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 12499.75
$ws.Range("J41").Value = 12499.75
$ws.Range("L41").Value = 12499.75
$ws.Range("N41").Value = -13279.75
$ws.Range("H108").Value = 56249.5
$ws.Range("J108").Value = 56249.5
$ws.Range("L108").Value = 56249.5
$ws.Range("N108").Value = -63929.5
$ws.Range("H122").Value = 112696.43
$ws.Range("I122").Value = 156775
$ws.Range("K122").Value = 470325
$ws.Range("M122").Value = -467875
